$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E1) is no longer produced by the scraper's
# xpath expressions, so remove it and shift the remaining header columns
# (reviews_average, latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one position to the left.
$ws.Range("E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
